# Trade #7 closed at 2026-02-17 20:48:19 - unknown UNKNOWN +0.000%
#
# A MarketMaking trade (row 36 on "All Trades" / row 3 on "MarketMaking")
# that was OPEN is now closed with an early exit, and a brand-new
# MarketMaking trade (row 69 on "All Trades" / row 36 on "MarketMaking")
# has been opened. The Summary and Strategy Status sheets are updated to
# reflect the new trade counts / win rate.

$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 35   # Total Trades
$wsSummary.Range("B9").Value = 40   # Win Rate %

# --- Sheet: Strategy Status ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D5").Value = 2     # MarketMaking Trades

# --- Sheet: All Trades ---
$wsAll = $wb.Worksheets.Item("All Trades")

# Close out the existing open trade (Trade #35) recorded in row 36.
$wsAll.Range("G36").Value = 0.9           # Exit Price
$wsAll.Range("H36").Value = "CLOSED"      # Status
$wsAll.Range("K36").Value = 100.3         # Capital After
$wsAll.Range("L36").Value = "early_exit"  # Exit Reason
$wsAll.Range("M36").Value = 0.13          # Duration (min)

# Append the brand-new trade (Trade #68) as row 69.
$wsAll.Range("A69").Value = 68
$wsAll.Range("B69").Value = "'2026-02-17"
$wsAll.Range("C69").Value = "20:48:12"
$wsAll.Range("D69").Value = "MarketMaking"
$wsAll.Range("E69").Value = "DOWN"
$wsAll.Range("F69").Value = 0.9
$wsAll.Range("H69").Value = "OPEN"
$wsAll.Range("I69").Value = 0
$wsAll.Range("J69").Value = 0
$wsAll.Range("K69").Value = 100.3009090909091
$wsAll.Range("M69").Value = 0
$wsAll.Range("N69").Value = 0
$wsAll.Range("O69").Value = 0
$wsAll.Range("P69").Value = 0.6
$wsAll.Range("Q69").Value = "Normal spread capture: 19600 bps"

# --- Sheet: MarketMaking ---
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Close out the existing open trade (Trade #35) recorded in row 3.
$wsMM.Range("G3").Value = 0.9             # Exit Price
$wsMM.Range("H3").Value = "CLOSED"        # Status
$wsMM.Range("K3").Value = 100.3           # Capital After
$wsMM.Range("P3").Value = "early_exit"    # Exit Reason
$wsMM.Range("Q3").Value = 0.13            # Duration (min)

# Append the brand-new trade (Trade #68) as row 36.
$wsMM.Range("A36").Value = 68
$wsMM.Range("B36").Value = "'2026-02-17"
$wsMM.Range("C36").Value = "20:48:12"
$wsMM.Range("D36").Value = "MarketMaking"
$wsMM.Range("E36").Value = "DOWN"
$wsMM.Range("F36").Value = 0.9
$wsMM.Range("H36").Value = "OPEN"
$wsMM.Range("I36").Value = 0
$wsMM.Range("J36").Value = 0
$wsMM.Range("K36").Value = 100.3009090909091
$wsMM.Range("L36").Value = 0
$wsMM.Range("M36").Value = 0
$wsMM.Range("N36").Value = 0.6
$wsMM.Range("O36").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("Q36").Value = 0
